{"js": "// Update team-meetings related hours, per the commit:\n//   \"Update team meetings tasks to reflect other status reports\"\n//\n//   1) Week-1 \"Total Hours\" run \"3h*\"  -> \"1h*\"   (total becomes \"21h*\")\n//   2) Task 5 \"Team meetings\" run \"5h\" -> \"4h\"\n//   3) Week-2 \"Total Hours\" run \"23h\"  -> \"22h\"\n//\n// Each change is, at the character level, a single-digit substitution\n// inside an existing (underlined) run, so every edit below narrows a\n// `search()` hit down to just the digit being changed and replaces only\n// that digit - this preserves the run's other characters/formatting\n// untouched, matching the semantic content of the original edit.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Hunk 1: \"3h*\" -> \"1h*\"  (Week 1 total hours: \"2\" + \"3h*\" => \"2\" + \"1h*\")\n// ---------------------------------------------------------------------\nconst hunk1 = body.search(\"3h*\", { matchCase: true });\nhunk1.load(\"items\");\nawait context.sync();\nif (hunk1.items.length !== 1) {\n  throw new Error(\"Expected exactly one '3h*' match, found \" + hunk1.items.length);\n}\nconst hunk1Digit = hunk1.items[0].search(\"3\", { matchCase: true });\nhunk1Digit.load(\"items\");\nawait context.sync();\nhunk1Digit.items[0].insertText(\"1\", \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Hunk 2: \"5h\" -> \"4h\" on the \"Team meetings\" task line (Task 5, week 2)\n// ---------------------------------------------------------------------\nconst hunk2 = body.search(\"5h\", { matchCase: true });\nhunk2.load(\"items\");\nawait context.sync();\n\nconst hunk2Paragraphs = hunk2.items.map((item) => item.paragraphs.getFirst());\nhunk2Paragraphs.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet teamMeetingsRange = null;\nfor (let i = 0; i < hunk2.items.length; i++) {\n  if (hunk2Paragraphs[i].text.indexOf(\"Team meetings\") !== -1) {\n    teamMeetingsRange = hunk2.items[i];\n    break;\n  }\n}\nif (!teamMeetingsRange) {\n  throw new Error(\"Could not find the 'Team meetings' task's '5h' hours run\");\n}\nconst hunk2Digit = teamMeetingsRange.search(\"5\", { matchCase: true });\nhunk2Digit.load(\"items\");\nawait context.sync();\nhunk2Digit.items[0].insertText(\"4\", \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Hunk 3: \"23h\" -> \"22h\"  (Week 2 total hours line, not the \"*\" one)\n// ---------------------------------------------------------------------\nconst hunk3 = body.search(\"23h\", { matchCase: true });\nhunk3.load(\"items\");\nawait context.sync();\n\nconst hunk3Paragraphs = hunk3.items.map((item) => item.paragraphs.getFirst());\nhunk3Paragraphs.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet weekTwoTotalRange = null;\nfor (let i = 0; i < hunk3.items.length; i++) {\n  if (hunk3Paragraphs[i].text.indexOf(\"*Additional tasks\") === -1) {\n    weekTwoTotalRange = hunk3.items[i];\n    break;\n  }\n}\nif (!weekTwoTotalRange) {\n  throw new Error(\"Could not find the week-2 'Total Hours' run '23h'\");\n}\nconst hunk3Digit = weekTwoTotalRange.search(\"3\", { matchCase: true });\nhunk3Digit.load(\"items\");\nawait context.sync();\nhunk3Digit.items[0].insertText(\"2\", \"Replace\");\nawait context.sync();\n", "ps1": "# Update team-meetings related hours, per the commit:\n#   \"Update team meetings tasks to reflect other status reports\"\n#\n#   1) Week-1 \"Total Hours\" run \"3h*\"  -> \"1h*\"   (total becomes \"21h*\")\n#   2) Task 5 \"Team meetings\" run \"5h\" -> \"4h\"\n#   3) Week-2 \"Total Hours\" run \"23h\"  -> \"22h\"\n#\n# Each change is, at the character level, a single-digit substitution\n# inside an existing (underlined) run, so every edit below narrows the\n# Find match down to just the digit being changed and replaces only that\n# digit - this preserves the run's other characters/formatting untouched.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($rng) {\n    $p = $rng.Duplicate\n    $p.Expand(4) | Out-Null   # wdParagraph\n    return $p.Text\n}\n\nfunction Replace-SingleDigit($matchRange, $oldDigit, $newDigit) {\n    $sub = $matchRange.Duplicate\n    $subFind = $sub.Find\n    $subFind.ClearFormatting()\n    $ok = $subFind.Execute($oldDigit)\n    if (-not $ok) { throw \"Could not locate digit '$oldDigit' inside matched range '$($matchRange.Text)'\" }\n    $sub.Text = $newDigit\n}\n\n# ---------------------------------------------------------------------\n# Hunk 1: \"3h*\" -> \"1h*\"  (Week 1 total hours: \"2\" + \"3h*\" => \"2\" + \"1h*\")\n# ---------------------------------------------------------------------\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$found1 = $find1.Execute(\"3h*\")\nif (-not $found1) { throw \"hunk1: '3h*' not found\" }\nReplace-SingleDigit $rng1 \"3\" \"1\"\n\n# ---------------------------------------------------------------------\n# Hunk 2: \"5h\" -> \"4h\" on the \"Team meetings\" task line (Task 5, week 2)\n# ---------------------------------------------------------------------\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$foundTeamMeetings = $false\nwhile ($find2.Execute(\"5h\")) {\n    $paraText = Get-ParaText $rng2\n    if ($paraText -like \"*Team meetings*\") {\n        Replace-SingleDigit $rng2 \"5\" \"4\"\n        $foundTeamMeetings = $true\n        break\n    }\n}\nif (-not $foundTeamMeetings) { throw \"hunk2: 'Team meetings' task's '5h' not found\" }\n\n# ---------------------------------------------------------------------\n# Hunk 3: \"23h\" -> \"22h\"  (Week 2 total hours line, not the \"*\" one)\n# ---------------------------------------------------------------------\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.ClearFormatting()\n$foundWeekTwoTotal = $false\nwhile ($find3.Execute(\"23h\")) {\n    $paraText = Get-ParaText $rng3\n    if ($paraText -notlike \"*Additional tasks*\") {\n        Replace-SingleDigit $rng3 \"3\" \"2\"\n        $foundWeekTwoTotal = $true\n        break\n    }\n}\nif (-not $foundWeekTwoTotal) { throw \"hunk3: week-2 'Total Hours' run '23h' not found\" }\n"}
